$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-28 Monday" "2024-10-29 Tuesday"

Replace-Text "828÷3=" "615÷6="
Replace-Text "401÷8=" "226÷2="
Replace-Text "191÷7=" "976÷3="
Replace-Text "296÷8=" "720÷7="
Replace-Text "277÷4=" "657÷4="
Replace-Text "223÷5=" "432÷5="
Replace-Text "123÷3=" "347÷7="
Replace-Text "564÷9=" "762÷3="
Replace-Text "320÷8=" "219÷6="
Replace-Text "662÷9=" "531÷7="
Replace-Text "586÷4=" "841÷8="
Replace-Text "506÷3=" "612÷3="
Replace-Text "833÷6=" "876÷4="
Replace-Text "887÷4=" "835÷9="
Replace-Text "740÷6=" "986÷4="
Replace-Text "822÷4=" "951÷2="
Replace-Text "917÷4=" "626÷4="
Replace-Text "324÷2=" "364÷3="
Replace-Text "231÷4=" "319÷9="
Replace-Text "995÷5=" "170÷4="
Replace-Text "415÷6=" "119÷7="
Replace-Text "515÷3=" "845÷9="
Replace-Text "470÷9=" "220÷7="
Replace-Text "844÷6=" "488÷5="
Replace-Text "354÷6=" "760÷6="

Write-Output "Done"
